$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 11.14494766666667
$ws.Range("H2").Value = 33.434843
$ws.Range("I2").Value = 0.1279818847384872
$ws.Range("J2").Value = 0.1279818847384872
$ws.Range("M2").Value = 0.952391
$ws.Range("N2").Value = 2.857173
$ws.Range("O2").Value = 0.1470615060759953
$ws.Range("P2").Value = 0.1470615060759953
$ws.Range("Q2").Value = 10.61434785320433
$ws.Range("R2").Value = 95.529130678839
$ws.Range("S2").Value = 0.01882120872008637
$ws.Range("T2").Value = 0.01882120872008637
$ws.Range("G3").Value = 11.14494766666667
$ws.Range("H3").Value = 33.434843
$ws.Range("I3").Value = 0.1279818847384872
$ws.Range("J3").Value = 0.1279818847384872
$ws.Range("O3").Value = 0.4906041777350729
$ws.Range("P3").Value = 0.4906041777350729
$ws.Range("Q3").Value = 35.40996920039944
$ws.Range("R3").Value = 318.689722803595
$ws.Range("S3").Value = 0.0627884473271104
$ws.Range("T3").Value = 0.0627884473271104
$ws.Range("G4").Value = 11.14494766666667
$ws.Range("H4").Value = 33.434843
$ws.Range("I4").Value = 0.1279818847384872
$ws.Range("J4").Value = 0.1279818847384872
$ws.Range("M4").Value = 2.346528
$ws.Range("N4").Value = 7.039584
$ws.Range("O4").Value = 0.3623343161889319
$ws.Range("P4").Value = 0.3623343161889319
$ws.Range("Q4").Value = 26.151931758368
$ws.Range("R4").Value = 235.367385825312
$ws.Range("S4").Value = 0.04637222869129046
$ws.Range("T4").Value = 0.04637222869129047
$ws.Range("I5").Value = 0.5307607770439682
$ws.Range("J5").Value = 0.5307607770439681
$ws.Range("M5").Value = 0.952391
$ws.Range("N5").Value = 2.857173
$ws.Range("O5").Value = 0.1470615060759953
$ws.Range("P5").Value = 0.1470615060759953
$ws.Range("Q5").Value = 44.019351065921
$ws.Range("R5").Value = 396.174159593289
$ws.Range("S5").Value = 0.0780544792381515
$ws.Range("T5").Value = 0.07805447923815148
$ws.Range("I6").Value = 0.5307607770439682
$ws.Range("J6").Value = 0.5307607770439681
$ws.Range("O6").Value = 0.4906041777350729
$ws.Range("P6").Value = 0.4906041777350729
$ws.Range("S6").Value = 0.2603934545956844
$ws.Range("T6").Value = 0.2603934545956843
$ws.Range("I7").Value = 0.5307607770439682
$ws.Range("J7").Value = 0.5307607770439681
$ws.Range("M7").Value = 2.346528
$ws.Range("N7").Value = 7.039584
$ws.Range("O7").Value = 0.3623343161889319
$ws.Range("P7").Value = 0.3623343161889319
$ws.Range("Q7").Value = 108.456127596768
$ws.Range("R7").Value = 976.1051483709119
$ws.Range("S7").Value = 0.1923128432101323
$ws.Range("T7").Value = 0.1923128432101323
$ws.Range("G8").Value = 29.71744933333333
$ws.Range("H8").Value = 89.152348
$ws.Range("I8").Value = 0.3412573382175446
$ws.Range("J8").Value = 0.3412573382175446
$ws.Range("M8").Value = 0.952391
$ws.Range("N8").Value = 2.857173
$ws.Range("O8").Value = 0.1470615060759953
$ws.Range("P8").Value = 0.1470615060759953
$ws.Range("Q8").Value = 28.30263128802267
$ws.Range("R8").Value = 254.723681592204
$ws.Range("S8").Value = 0.05018581811775742
$ws.Range("T8").Value = 0.05018581811775742
$ws.Range("G9").Value = 29.71744933333333
$ws.Range("H9").Value = 89.152348
$ws.Range("I9").Value = 0.3412573382175446
$ws.Range("J9").Value = 0.3412573382175446
$ws.Range("O9").Value = 0.4906041777350729
$ws.Range("P9").Value = 0.4906041777350729
$ws.Range("Q9").Value = 94.41892389993556
$ws.Range("R9").Value = 849.77031509942
$ws.Range("S9").Value = 0.1674222758122781
$ws.Range("T9").Value = 0.1674222758122781
$ws.Range("G10").Value = 29.71744933333333
$ws.Range("H10").Value = 89.152348
$ws.Range("I10").Value = 0.3412573382175446
$ws.Range("J10").Value = 0.3412573382175446
$ws.Range("M10").Value = 2.346528
$ws.Range("N10").Value = 7.039584
$ws.Range("O10").Value = 0.3623343161889319
$ws.Range("P10").Value = 0.3623343161889319
$ws.Range("Q10").Value = 69.732826949248
$ws.Range("R10").Value = 627.595442543232
$ws.Range("S10").Value = 0.1236492442875091
$ws.Range("T10").Value = 0.1236492442875091
